$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for a handful of rows, as part of a
# repull/recalculation of the underlying data.
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = 7
